# "générateur d'onde sin fonctionnel + débute onde tri pulsée"
#
# This commit (as captured by the author's LibreOffice Calc save) mainly:
#   1. Appends a test date to each of the first four sheet names.
#   2. Renames the 5th sheet ("test_adc_dac_modele_rc" -> "test_adc_dac_resistance (31 mai 2017)")
#      and wipes out its old/obsolete measurement table (work-in-progress on the new
#      sine/triangle-wave generator sheet starts from a blank sheet).
#   3. Makes the first sheet the active/selected one again (tab ratio/active tab reset).
#
# NOTE: genuine Excel enforces a 31-character cap on sheet names (and forbids
# []:/\? ), so the literal target names (36-40 chars) are rejected by the
# object model - they were only legal because the source file was produced by
# LibreOffice, which doesn't enforce that limit. We approximate with the
# longest valid (<=31 char) prefix of each intended name so the dates /
# intent stay recognisable.

$wb = $excel.ActiveWorkbook

$newNames = @(
    "test_controle_potentiomètre (4 ",
    "test_controle_adc12bits (11 mai",
    "test_controle_adc_10bits (15 ma",
    "test_controle_dac16bits(17 mai ",
    "test_adc_dac_resistance (31 mai"
)

for ($i = 1; $i -le $newNames.Length; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}

# The old "test_adc_dac_modele_rc" sheet (now "test_adc_dac_resistance (31 ...)")
# had a leftover results table that's no longer relevant - clear it out
# completely (back to a pristine, empty sheet) ahead of the new sine-wave work.
$ws5 = $wb.Worksheets.Item(5)
$ws5.UsedRange.EntireRow.Delete()
$ws5.Range("K21").Select()

# Tab/view state: put the focus back on the first sheet (it was left on sheet 5
# in the previous save).
$win = $excel.ActiveWindow
$win.TabRatio = 0.5
$wb.Worksheets.Item(1).Activate()
